{"js": "// Replace the recorrente's surname \"CARMO\" with \"CAMARGO\" everywhere it occurs\n// (matching the exact full name to avoid touching unrelated text).\nconst nameResults = context.document.body.search(\n  \"GUILHERME J\u00daNIOR DA SILVA CARMO\",\n  { matchCase: true }\n);\nnameResults.load(\"text\");\nawait context.sync();\n\nnameResults.items.forEach((r) => {\n  r.insertText(\"GUILHERME J\u00daNIOR DA SILVA CAMARGO\", Word.InsertLocation.replace);\n});\nawait context.sync();\n\n// Replace the bullet-list paragraphs under \"DO M\u00c9RITO RECURSAL\" with the new\n// argument summary text. Map old (unique) paragraph text -> new paragraph text.\nconst bulletReplacements = [\n  [\n    \"- Aus\u00eancia de prequestionamento da tese de desclassifica\u00e7\u00e3o para recepta\u00e7\u00e3o culposa.\",\n    \"- Aus\u00eancia de prequestionamento da mat\u00e9ria relativa \u00e0 recepta\u00e7\u00e3o culposa (art. 180, \u00a7 3\u00ba, CP)\",\n  ],\n  [\n    \"- Incid\u00eancia da S\u00famula 7 do STJ quanto \u00e0 an\u00e1lise da ilicitude da prova.\",\n    \"- Incid\u00eancia da S\u00famula 7 do STJ quanto \u00e0 an\u00e1lise da licitude da busca domiciliar\",\n  ],\n  [\n    \"- Incid\u00eancia da S\u00famula 7 do STJ quanto \u00e0 an\u00e1lise da habitualidade da atividade comercial.\",\n    \"- Incid\u00eancia da S\u00famula 7 do STJ quanto \u00e0 caracteriza\u00e7\u00e3o da atividade comercial\",\n  ],\n  [\n    \"- Inexist\u00eancia de flagrante ilegalidade na busca domiciliar.\",\n    \"- Incid\u00eancia da S\u00famula 7 do STJ quanto ao conhecimento da origem il\u00edcita do bem\",\n  ],\n  [\n    \"- Decis\u00e3o recorrida em conson\u00e2ncia com o entendimento do STJ sobre a recepta\u00e7\u00e3o qualificada.\",\n    \"- Inexist\u00eancia de viola\u00e7\u00e3o ao art. 157 do CPP: legalidade da busca domiciliar com base em fundada suspeita\",\n  ],\n  [\n    \"- Reexame de provas invi\u00e1vel em Recurso Especial.\",\n    \"- Inexist\u00eancia de viola\u00e7\u00e3o ao art. 180, \u00a7 1\u00ba, do CP: correta caracteriza\u00e7\u00e3o da atividade comercial\",\n  ],\n  [\n    \"- S\u00famula 83/STJ.\",\n    \"- Inexist\u00eancia de viola\u00e7\u00e3o ao art. 180, caput, do CP: demonstra\u00e7\u00e3o do dolo na conduta do agente\",\n  ],\n  [\n    \"- M\u00e9rito: Validade da busca domiciliar diante das fundadas suspeitas.\",\n    \"- Preval\u00eancia do princ\u00edpio do livre convencimento motivado do juiz\",\n  ],\n  [\n    \"- M\u00e9rito: Sufici\u00eancia de provas para a condena\u00e7\u00e3o por recepta\u00e7\u00e3o qualificada.\",\n    \"- M\u00e9rito: Sufici\u00eancia das provas para a condena\u00e7\u00e3o por recepta\u00e7\u00e3o qualificada\",\n  ],\n  [\n    \"- M\u00e9rito: Dolo comprovado na conduta do r\u00e9u.\",\n    \"- M\u00e9rito: Aplica\u00e7\u00e3o do princ\u00edpio *pas de nullit\u00e9 sans grief*\",\n  ],\n];\n\nconst paragraphs = context.document.body.paragraphs;\nparagraphs.load(\"text\");\nawait context.sync();\n\nfor (const para of paragraphs.items) {\n  for (const [oldText, newText] of bulletReplacements) {\n    if (para.text === oldText) {\n      para.insertText(newText, Word.InsertLocation.replace);\n      break;\n    }\n  }\n}\nawait context.sync();\n", "ps1": "$d = $word.ActiveDocument\n\nfunction Replace-Text($findText, $replaceText) {\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Text = $findText\n    $find.Replacement.Text = $replaceText\n    # wdFindContinue=1, wdReplaceAll=2\n    $find.Execute($null, $true, $false, $false, $false, $false, $true, 1, $false, $null, 2)\n}\n\n# Surname correction: CARMO -> CAMARGO (applies to every occurrence of the\n# recorrente's full name throughout the document).\nReplace-Text \"GUILHERME J\u00daNIOR DA SILVA CARMO\" \"GUILHERME J\u00daNIOR DA SILVA CAMARGO\"\n\n# Updated bullet list under \"DO M\u00c9RITO RECURSAL\" with the new argument summary.\nReplace-Text \"- Aus\u00eancia de prequestionamento da tese de desclassifica\u00e7\u00e3o para recepta\u00e7\u00e3o culposa.\" \"- Aus\u00eancia de prequestionamento da mat\u00e9ria relativa \u00e0 recepta\u00e7\u00e3o culposa (art. 180, \u00a7 3\u00ba, CP)\"\nReplace-Text \"- Incid\u00eancia da S\u00famula 7 do STJ quanto \u00e0 an\u00e1lise da ilicitude da prova.\" \"- Incid\u00eancia da S\u00famula 7 do STJ quanto \u00e0 an\u00e1lise da licitude da busca domiciliar\"\nReplace-Text \"- Incid\u00eancia da S\u00famula 7 do STJ quanto \u00e0 an\u00e1lise da habitualidade da atividade comercial.\" \"- Incid\u00eancia da S\u00famula 7 do STJ quanto \u00e0 caracteriza\u00e7\u00e3o da atividade comercial\"\nReplace-Text \"- Inexist\u00eancia de flagrante ilegalidade na busca domiciliar.\" \"- Incid\u00eancia da S\u00famula 7 do STJ quanto ao conhecimento da origem il\u00edcita do bem\"\nReplace-Text \"- Decis\u00e3o recorrida em conson\u00e2ncia com o entendimento do STJ sobre a recepta\u00e7\u00e3o qualificada.\" \"- Inexist\u00eancia de viola\u00e7\u00e3o ao art. 157 do CPP: legalidade da busca domiciliar com base em fundada suspeita\"\nReplace-Text \"- Reexame de provas invi\u00e1vel em Recurso Especial.\" \"- Inexist\u00eancia de viola\u00e7\u00e3o ao art. 180, \u00a7 1\u00ba, do CP: correta caracteriza\u00e7\u00e3o da atividade comercial\"\nReplace-Text \"- S\u00famula 83/STJ.\" \"- Inexist\u00eancia de viola\u00e7\u00e3o ao art. 180, caput, do CP: demonstra\u00e7\u00e3o do dolo na conduta do agente\"\nReplace-Text \"- M\u00e9rito: Validade da busca domiciliar diante das fundadas suspeitas.\" \"- Preval\u00eancia do princ\u00edpio do livre convencimento motivado do juiz\"\nReplace-Text \"- M\u00e9rito: Sufici\u00eancia de provas para a condena\u00e7\u00e3o por recepta\u00e7\u00e3o qualificada.\" \"- M\u00e9rito: Sufici\u00eancia das provas para a condena\u00e7\u00e3o por recepta\u00e7\u00e3o qualificada\"\nReplace-Text \"- M\u00e9rito: Dolo comprovado na conduta do r\u00e9u.\" \"- M\u00e9rito: Aplica\u00e7\u00e3o do princ\u00edpio *pas de nullit\u00e9 sans grief*\"\n"}
